$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# The original sheet has a single merged column-width definition covering
# A:C ("25.7109375"), even though only A:B actually hold data (dimension
# is A1:B99). The edit widens column A (to fit the long "DisplayName"
# values) and drops the unused, empty column C definition entirely.
$ws.Columns("C").Delete()
$ws.Columns("A").ColumnWidth = 79.83333333333333
$ws.Columns("B").ColumnWidth = 24.833333333333332

# --- Header formatting ---------------------------------------------------
# Strip the manual bold+boxed-border formatting that was applied directly
# to the header cells; the table style (added below) takes over the
# header look instead.
$ws.Range("A1:B1").ClearFormats()

# --- Convert the data range into a real Excel Table ("ListObject") -----
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:B99"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "TabelaSoftware"
$tbl.TableStyle = "TableStyleMedium9"

# --- Freeze the header rows --------------------------------------------
# Freeze the top two rows (DisplayName/DisplayVersion header + the
# "Usuário Logado" info row), leaving row 3 as the first scrollable row.
$ws.Activate()
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true
